$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 617 (shifts the existing rows 617..650 down to 618..651)
$ws.Rows(617).Insert()

$r = 617
$ws.Cells.Item($r, 1).Value  = 9
$ws.Cells.Item($r, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($r, 3).Value  = "Metropolitana"
$ws.Cells.Item($r, 4).Value  = 45041
$ws.Cells.Item($r, 5).Value  = 13
$ws.Cells.Item($r, 6).Value  = "Fruta"
$ws.Cells.Item($r, 7).Value  = 100108
$ws.Cells.Item($r, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value  = 100108002
$ws.Cells.Item($r, 10).Value = "Mango"
$ws.Cells.Item($r, 11).Value = "Sin especificar"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 600
$ws.Cells.Item($r, 14).Value = 7000
$ws.Cells.Item($r, 15).Value = 7500
$ws.Cells.Item($r, 16).Value = 7250
$ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item($r, 18).Value = "Perú"
$ws.Cells.Item($r, 19).Value = 1812
$ws.Cells.Item($r, 20).Value = 4
